# Generate Report for handoff
# The d189dbdd-...md file is now ready for handoff: update the Status
# columns on the Overview / zh-cn / de-de sheets, and refresh the
# "Latest Handoff Datetime" stamps for the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for d189dbdd-...md is row 3 ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 3 Status + Latest Handoff Datetime ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-01-28 09:13:16"

# --- de-de sheet: row 3 Status + Latest Handoff Datetime ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-01-28 09:13:29"
